$d = $word.ActiveDocument

# --- Paragraph: "Put slice of bread with peanut butter against slice of ..." ---
$xmlPut = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="46"/>
    </w:numPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Put </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t>slice of bread with peanut butter against slice of bread where</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t xml:space="preserve"> peanut butter and </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t>other side of bread</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t xml:space="preserve"> are touching so that you maximize the area of contact</w:t>
  </w:r>
</w:p>
'@

# --- Paragraph: "Pre Condition: Have prepared slices of bread with peanut butter and ..." ---
$xmlPre = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="46"/>
    </w:numPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t>Pre Condition: Have prepared slic</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t xml:space="preserve">es of bread with peanut butter and </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t>the other plain slice of bread</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
'@

# --- Paragraph: "Post Condition: Have a full, eatable sandwich" ---
$xmlPost = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="46"/>
    </w:numPr>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t>Post Condition: Have a full,</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Helvetica Neue"/>
    </w:rPr>
    <w:t xml:space="preserve"> eatable sandwich</w:t>
  </w:r>
</w:p>
'@

# --- Paragraph: "Cut sandwich in four different pieces in a cross manner" (now bold) ---
$xmlCut = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="46"/>
    </w:numPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:b/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
    </w:rPr>
    <w:t>Cut sandwich in four different pieces in a cross manner</w:t>
  </w:r>
</w:p>
'@

# --- Big closing paragraph ("As before, ...") + the two trailing empty paragraphs ---
$xmlClosing = @'
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>As before</w:t>
  </w:r>
  <w:r>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> we decided to show the updates by bolding the lettering. As far as our major changes</w:t>
  </w:r>
  <w:r>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> we saved a bunch of time when it came to the spreading of the peanut butter.  We decided to use slices of peanut butter in order </w:t>
  </w:r>
  <w:r>
    <w:t>to cut</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> the time it would take </w:t>
  </w:r>
  <w:r>
    <w:t>to spread</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">  This also ensures quality and consistency for each sandwich.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">  Furthermore</w:t>
  </w:r>
  <w:r>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> we decided to cut the sandwich into four pieces and only need ¼ of the sandwich in order to test</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> using a special pizza cutting utensil that cuts down our slicing time as well as providing consistent quality with every cut</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t xml:space="preserve">. </w:t>
  </w:r>
  <w:r>
    <w:t>We can then distribute the rest of the sandwich to our customers</w:t>
  </w:r>
  <w:r>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> which was a problem before.  </w:t>
  </w:r>
  <w:r>
    <w:t>Even we have changed the pr</w:t>
  </w:r>
  <w:r>
    <w:t>ocess in order to save time as well as outsourcing to countries that cost ¾ less than local manufactures,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> we</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> were still over budget. </w:t>
  </w:r>
  <w:r>
    <w:t>The outsourcing significantly decreased our overhead,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> but still didn’t get us within budget. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> The real issue with the offshore team </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">is </w:t>
  </w:r>
  <w:r>
    <w:t>having 3x quality issues</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">.  We would need to add another team offshore in order to account for the 30% loss in product.  In the end, this could add up to $1 million dollars extra (the entirety of our budget).  </w:t>
  </w:r>
  <w:r>
    <w:t>Even though the offshore team is 25% cost compared to the local team, t</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">hat 30% error is </w:t>
  </w:r>
  <w:r>
    <w:t>a detrimental hit to our development</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> and budget</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
  </w:pPr>
</w:p>
'@

function Replace-ParagraphXml($paragraph, $xml) {
    $paragraph.Range.InsertXML($xml)
}

# The "Put slice of bread..." / "Pre Condition: ... jelly" / "Post Condition:
# Have an eatable sandwich" / "Cut sandwich in four different pieces..."
# wording recurs several times in this document (earlier steps of the same
# recipe reuse near-identical boilerplate). Only the final occurrence - the
# one that already mentions "jelly" - is the one touched by this edit, so we
# anchor on that unique "jelly" paragraph and then walk forward from it
# (the three paragraphs that follow it, in order, are the Post Condition /
# Cut sandwich ones we need) instead of doing independent whole-document
# text searches that would also hit the earlier, unrelated occurrences.

$putIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -match "^Put slice of bread with peanut butter against.*jelly where peanut butter and jelly are touching") {
        $putIndex = $i
        break
    }
}

if ($putIndex -eq -1) {
    throw "Could not locate the 'Put slice of bread ... jelly ...' paragraph"
}

Replace-ParagraphXml $d.Paragraphs($putIndex) $xmlPut
Replace-ParagraphXml $d.Paragraphs($putIndex + 1) $xmlPre
Replace-ParagraphXml $d.Paragraphs($putIndex + 2) $xmlPost
Replace-ParagraphXml $d.Paragraphs($putIndex + 3) $xmlCut

# The closing "As before, ..." paragraph is followed immediately by the two
# trailing blank paragraphs (the second of which used to hold the _GoBack
# bookmark). The edit folds that bookmark into the middle of the "As
# before..." paragraph's runs, so all three paragraphs must be replaced as a
# single block (1 paragraph growing into 3) rather than touched one at a
# time - otherwise the original trailing two paragraphs (and their
# now-duplicated bookmark) would be left behind.
$asBeforeIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -match "^As before, we decided to show the updates by bolding the lettering") {
        $asBeforeIndex = $i
        break
    }
}

if ($asBeforeIndex -eq -1) {
    throw "Could not locate the 'As before, ...' paragraph"
}

$firstPara = $d.Paragraphs($asBeforeIndex)
$lastPara = $d.Paragraphs($asBeforeIndex + 2)
$blockRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
$blockRange.InsertXML($xmlClosing)

Write-Host "Done"
